$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "fix accessions with a dash for a range of values" - five proposal
# packages whose ACCESSION field used a dash to denote a range; add one
# row per affected package documenting the fix (columns B..G), and extend
# column A's "All <X> proposals" helper formula down to match.
#
# Cell writes are ordered to match the authoring sequence (new shared
# strings get interned in this order: B20, C20, B21, G20, B22, B23, B24).
$ws.Cells.Item(20, 2).Value = "2021.004P.A.v1.Anulavirus_1ns.zip"
$ws.Cells.Item(20, 3).Value = "ACCESSION range"
$ws.Cells.Item(21, 2).Value = "2021.005P.A.v1.Bromovirus_1ns.zip"
$ws.Cells.Item(20, 7).Value = "fix3"
$ws.Cells.Item(22, 2).Value = "2021.010P.A.v1.Emaravirus_1ns.zip"
$ws.Cells.Item(23, 2).Value = "2021.011P.A.v1.Emaravirus_1ns.zip"
$ws.Cells.Item(24, 2).Value = "2021.012P.A.v1.Emaravirus_1ns.zip"

for ($r = 20; $r -le 24; $r++) {
    $ws.Cells.Item($r, 3).Value = "ACCESSION range"
    $ws.Cells.Item($r, 4).Value = "yes"
    $ws.Cells.Item($r, 5).Value = "yes"
    $ws.Cells.Item($r, 6).Value = "yes"
    $ws.Cells.Item($r, 7).Value = "fix3"
}

# Extend the shared "All X proposals" helper formula in column A down
# through the newly added rows (previously A3:A19, now A3:A24).
$ws.Range("A20:A24").Formula = '=CONCATENATE("All ",MID(B20,9,1)," proposals")'

# Leave the selection on the blank row right after the new data, like a
# user would after finishing data entry.
$ws.Range("A25:XFD25").Select() | Out-Null
